# "updates after talking to Patrick"
#
# Adds 5 new rows to the PLLF tracker table (Table1) and updates the
# "action" text + row height of the last existing row (row 25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

function Add-TrackerRow {
    param($Problem, $RaisedBy, $NotedDate, $Action, $RowHeight)

    $newRow = $lo.ListRows.Add()
    $rowNum = $newRow.Range.Row

    # Copy formatting (wrap text / top align / date format) from the
    # previous last data row (row 25) so the new row's styles match
    # the existing table styling instead of Excel's defaults.
    $srcRow = $ws.Range("A25:E25")
    $dstRow = $ws.Range("A" + $rowNum + ":E" + $rowNum)
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("A" + $rowNum).Value = $Problem
    $ws.Range("B" + $rowNum).Value = $RaisedBy
    $ws.Range("C" + $rowNum).Value = $NotedDate
    $ws.Range("D" + $rowNum).Value = $Action

    if ($RowHeight) {
        $ws.Rows($rowNum).RowHeight = $RowHeight
    }

    return $rowNum
}

# Row 26: "save the normal approx to data?"
$dummy = Add-TrackerRow "save the normal approx to data?" "Ian" 45874 "discuss" $null

# Row 25 (existing last row): action text changes + taller row
$ws.Range("D25").Value = "Reload the model with the selected parameter value (MLE) and leave that model in memory.  If there is no such value we just leave nothing via ereturn clear. "
$ws.Rows("25").RowHeight = 90

# Row 27
$dummy = Add-TrackerRow "change first help file example from 38-19 to 8-3" "Ian" 45883 "" $null

# Row 28
$dummy = Add-TrackerRow 'Improve "Stored results" section - to be like -h regress-? ' "Ian" 45883 "" 30

# Row 29
$dummy = Add-TrackerRow "Improve help file by showing both syntaxes in the syntax diagram, thus also explaining what pllf does; changing the RHS of the syntax diagram to just regression_cmd" "Ian" 45883 "" 60

# Row 30
$dummy = Add-TrackerRow "Syntax 2: Remove the placeholder X from regression_cmd; make the default placeholder(#)" "Ian" 45883 "" 30

# Match the final selection/scroll state from the authored workbook.
$dummy = $ws.Range("B28:C30").Select()
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
